$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) from 45181 to 45182 for all data rows (2..308)
for ($r = 2; $r -le 308; $r++) {
    $ws.Cells.Item($r, 3).Value = 45182
}

# Row 308 gains an explicit row height (15, custom height)
$ws.Rows.Item(308).RowHeight = 15

# New row 309
$ws.Rows.Item(309).RowHeight = 15
$ws.Cells.Item(309, 1).Value = "A 42455-2023"
$ws.Cells.Item(309, 2).Value = 45180
$ws.Cells.Item(309, 2).NumberFormat = $ws.Cells.Item(308, 2).NumberFormat
$ws.Cells.Item(309, 3).Value = 45182
$ws.Cells.Item(309, 3).NumberFormat = $ws.Cells.Item(308, 3).NumberFormat
$ws.Cells.Item(309, 4).Value = "HALLANDS LÄN"
$ws.Cells.Item(309, 5).Value = "LAHOLM"
$ws.Cells.Item(309, 6).Value = "Sveaskog"
$ws.Cells.Item(309, 7).Value = 3.6
$ws.Cells.Item(309, 8).Value = 0
$ws.Cells.Item(309, 9).Value = 0
$ws.Cells.Item(309, 10).Value = 0
$ws.Cells.Item(309, 11).Value = 0
$ws.Cells.Item(309, 12).Value = 0
$ws.Cells.Item(309, 13).Value = 0
$ws.Cells.Item(309, 14).Value = 0
$ws.Cells.Item(309, 15).Value = 0
$ws.Cells.Item(309, 16).Value = 0
$ws.Cells.Item(309, 17).Value = 0
$ws.Cells.Item(309, 18).WrapText = $true

# New row 310
$ws.Cells.Item(310, 1).Value = "A 42605-2023"
$ws.Cells.Item(310, 2).Value = 45181
$ws.Cells.Item(310, 2).NumberFormat = $ws.Cells.Item(308, 2).NumberFormat
$ws.Cells.Item(310, 3).Value = 45182
$ws.Cells.Item(310, 3).NumberFormat = $ws.Cells.Item(308, 3).NumberFormat
$ws.Cells.Item(310, 4).Value = "HALLANDS LÄN"
$ws.Cells.Item(310, 5).Value = "LAHOLM"
$ws.Cells.Item(310, 7).Value = 0.7
$ws.Cells.Item(310, 8).Value = 0
$ws.Cells.Item(310, 9).Value = 0
$ws.Cells.Item(310, 10).Value = 0
$ws.Cells.Item(310, 11).Value = 0
$ws.Cells.Item(310, 12).Value = 0
$ws.Cells.Item(310, 13).Value = 0
$ws.Cells.Item(310, 14).Value = 0
$ws.Cells.Item(310, 15).Value = 0
$ws.Cells.Item(310, 16).Value = 0
$ws.Cells.Item(310, 17).Value = 0
$ws.Cells.Item(310, 18).WrapText = $true
